$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.912.37'
$ws.Range("E2").Value = '  -1.63%  '

# Row 3
$ws.Range("D3").Value = '2.543.69'
$ws.Range("E3").Value = '  +0.41%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.39'
$ws.Range("E5").Value = '  -0.61%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.00'
$ws.Range("E6").Value = '  +0.32%  '

# Row 7
$ws.Range("E7").Value = '  +0.08%  '

# Row 8
$ws.Range("E8").Value = '  -1.31%  '

# Row 9
$ws.Range("D9").Value = '2.543.02'
$ws.Range("E9").Value = '  +0.62%  '

# Row 10
$ws.Range("E10").Value = '  -1.49%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.49'
$ws.Range("E11").Value = '  -4.97%  '

# Row 12
$ws.Range("E12").Value = '  -0.43%  '

# Row 13
$ws.Range("E13").Value = '  -1.37%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.20'
$ws.Range("E14").Value = '  -2.73%  '

# Row 15
$ws.Range("D15").Value = '2.998.58'
$ws.Range("E15").Value = '  +0.74%  '

# Row 16
$ws.Range("D16").Value = '62.825.93'
$ws.Range("E16").Value = '  -1.35%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000142'
$ws.Range("E17").Value = '  -0.99%  '

# Row 18
$ws.Range("D18").Value = '2.545.22'
$ws.Range("E18").Value = '  +1.43%  '

# Row 19
$ws.Range("E19").Value = '  -1.59%  '

# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '334.14'
$ws.Range("E20").Value = '  -2.86%  '

# Row 21
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.32'
$ws.Range("E21").Value = '  -0.65%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.75'
$ws.Range("E22").Value = '  -2.15%  '

# Row 23
$ws.Range("E23").Value = '  -0.42%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.37'
$ws.Range("E24").Value = '  -0.99%  '

# Row 25
$ws.Range("E25").Value = '  -1.55%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.58'
$ws.Range("E26").Value = '  +2.91%  '

# Row 27
$ws.Range("E27").Value = '  +0.42%  '

# Row 28
$ws.Range("E28").Value = '  +1.27%  '

# Row 29
$ws.Range("E29").Value = '  +2.73%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.29'
$ws.Range("E30").Value = '  +7.41%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0809'
$ws.Range("E31").Value = '  -2.13%  '

# Row 32
$ws.Range("E32").Value = '  -1.83%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '175.67'
$ws.Range("E33").Value = '  -0.44%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.53'
$ws.Range("E34").Value = '  -0.42%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '404.87'
$ws.Range("E35").Value = '  -1.17%  '

# Row 36
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.07'
$ws.Range("E36").Value = '  -0.32%  '

# Row 37
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.399'
$ws.Range("E37").Value = '  -1.10%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.32'
$ws.Range("E39").Value = '  -2.37%  '

# Row 40
$ws.Range("E40").Value = '  -0.83%  '

# Row 41
$ws.Range("E41").Value = '  +0.16%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.53'
$ws.Range("E42").Value = '  -3.05%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '151.50'
$ws.Range("E43").Value = '  -0.94%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.73'
$ws.Range("E44").Value = '  -1.72%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.75'
$ws.Range("E45").Value = '  -1.06%  '

# Row 46
$ws.Range("E46").Value = '  +1.04%  '

# Row 47
$ws.Range("E47").Value = '  -1.77%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0965'
$ws.Range("E48").Value = '  -0.41%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0239'
$ws.Range("E49").Value = '  +3.11%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.19'
$ws.Range("E50").Value = '  -3.25%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.72'
$ws.Range("E51").Value = '  -4.65%  '

